# Pio's ERA operativo Abril-Diciembre 2025
#
# Refresh the "Antiguedad" (seniority, in years) and the derived
# "Carga Teorica" (theoretical workload) columns for every médico on the
# sheet, using the new reference date for the ERA (turno) calculation
# cycle. A couple of the monthly-rate helper columns (Mañanas/mes,
# Tardes/mes, Noches/mes) that depend on the updated workload are
# refreshed as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Antiguedad (column C) - years of seniority as of the new reference date
$ws.Range("C2").Value  = 10.92328767123288
$ws.Range("C3").Value  = 10.92328767123288
$ws.Range("C4").Value  = 10.41917808219178
$ws.Range("C5").Value  = 10.41917808219178
$ws.Range("C6").Value  = 9.838356164383562
$ws.Range("C7").Value  = 9.838356164383562
$ws.Range("C8").Value  = 9.838356164383562
$ws.Range("C9").Value  = 7.416438356164384
$ws.Range("C10").Value = 6.583561643835616
$ws.Range("C11").Value = 6.383561643835616
$ws.Range("C12").Value = 4.583561643835616
$ws.Range("C13").Value = 2.747945205479452
$ws.Range("C14").Value = 2
$ws.Range("C15").Value = 2
$ws.Range("C16").Value = 0.1643835616438356

# Carga Teorica (column F) - theoretical load recomputed from the new
# Antiguedad values
$ws.Range("F2").Value  = 6.224379325864565
$ws.Range("F3").Value  = 6.224379325864565
$ws.Range("F4").Value  = 6.280334269970441
$ws.Range("F5").Value  = 6.280334269970441
$ws.Range("F6").Value  = 6.344804096875039
$ws.Range("F7").Value  = 6.344804096875039
$ws.Range("F8").Value  = 6.344804096875039
$ws.Range("F9").Value  = 6.613631110948924
$ws.Range("F10").Value = 6.706078409906459
$ws.Range("F11").Value = 6.728277925774551
$ws.Range("F12").Value = 6.928073568587381
$ws.Range("F13").Value = 7.131822549842476
$ws.Range("F14").Value = 7.214842657129998
$ws.Range("F15").Value = 7.214842657129998
$ws.Range("F16").Value = 7.418591638385093

# Derived monthly-rate helper columns that shift along with the workload
$ws.Range("O2").Value  = 0.3415730337078651
$ws.Range("P10").Value = 0.2841121495327102
$ws.Range("R13").Value = 0.3304347826086956
